$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Attack (column C) stats per new balancing values
$ws.Range("C2").Value = 100
$ws.Range("C3").Value = 120
$ws.Range("C4").Value = 140
$ws.Range("C6").Value = 180
$ws.Range("C7").Value = 200
$ws.Range("C8").Value = 220
$ws.Range("C9").Value = 240
$ws.Range("C10").Value = 260
$ws.Range("C11").Value = 280

# Restore the selection left on the sheet when it was saved
$ws.Range("E7").Select()
